$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 139; this shifts existing rows 139-163 down to 140-164
$ws.Rows(139).Insert()

# Populate the newly inserted row 139 with the new week's data
$ws.Cells.Item(139, 1).Value = 10
$ws.Cells.Item(139, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(139, 3).Value = "La Araucanía"
$ws.Cells.Item(139, 4).Value = 44508
$ws.Cells.Item(139, 5).Value = 9
$ws.Cells.Item(139, 6).Value = "Fruta"
$ws.Cells.Item(139, 7).Value = 100102
$ws.Cells.Item(139, 8).Value = "Cítricos"
$ws.Cells.Item(139, 9).Value = 100102006
$ws.Cells.Item(139, 10).Value = "Pomelo"
$ws.Cells.Item(139, 11).Value = "Start Ruby"
$ws.Cells.Item(139, 12).Value = "Primera"
$ws.Cells.Item(139, 13).Value = 55
$ws.Cells.Item(139, 14).Value = 13000
$ws.Cells.Item(139, 15).Value = 14000
$ws.Cells.Item(139, 16).Value = 13364
$ws.Cells.Item(139, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(139, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(139, 19).Value = 891
$ws.Cells.Item(139, 20).Value = 15
